# Daily cryptos data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.031.84'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.538.76'
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '198.83'
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '557.09'
$ws.Range("E6").Value = '  -2.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.647'
$ws.Range("E7").Value = '  +4.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.532.39'
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.664'
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '61.56'
$ws.Range("E11").Value = '  +6.73%  '
$ws.Range("E12").Value = '  -6.09%  '
$ws.Range("E13").Value = '  -7.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.97'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.092.05'
$ws.Range("E15").Value = '  -3.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.529.29'
$ws.Range("E16").Value = '  -3.11%  '
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.721.20'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.47'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.94'
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("E21").Value = '  -4.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '398.96'
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("E23").Value = '  -4.82%  '
$ws.Range("E24").Value = '  -8.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.76'
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("E26").Value = '  -1.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.86'
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.85'
$ws.Range("E28").Value = '  -4.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.93'
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '727.26'
$ws.Range("E30").Value = '  +3.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.42'
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.11'
$ws.Range("E32").Value = '  -13.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.80'
$ws.Range("E33").Value = '  -3.76%  '
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  -4.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.81'
$ws.Range("E36").Value = '  -9.28%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -7.30%  '
$ws.Range("E39").Value = '  -6.39%  '
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.090.43'
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0692'
$ws.Range("E43").Value = '  -12.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.52'
$ws.Range("E44").Value = '  -10.67%  '
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  -2.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.134'
$ws.Range("E47").Value = '  +1.82%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("E48").Value = '  -14.79%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.67'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.99'
$ws.Range("E50").Value = '  -2.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.28'
$ws.Range("E51").Value = '  -7.38%  '
